$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$cell = $t.Cell(1, 1)
$old = $cell.Range.Text
if ($old.Substring(0, $old.Length - 2) -ne '88÷5=17, 3') {
    throw "Cell(1,1) expected '88÷5=17, 3' but found '$old'"
}
$cell.Range.Text = '64÷3=21, 1'

$cell = $t.Cell(1, 2)
$old = $cell.Range.Text
if ($old.Substring(0, $old.Length - 2) -ne '37÷2=18, 1') {
    throw "Cell(1,2) expected '37÷2=18, 1' but found '$old'"
}
$cell.Range.Text = '67÷8=8, 3'

$cell = $t.Cell(1, 3)
$old = $cell.Range.Text
if ($old.Substring(0, $old.Length - 2) -ne '88÷4=22, 0') {
    throw "Cell(1,3) expected '88÷4=22, 0' but found '$old'"
}
$cell.Range.Text = '75÷7=10, 5'

$cell = $t.Cell(1, 4)
$old = $cell.Range.Text
if ($old.Substring(0, $old.Length - 2) -ne '56÷3=18, 2') {
    throw "Cell(1,4) expected '56÷3=18, 2' but found '$old'"
}
$cell.Range.Text = '42÷4=10, 2'

$cell = $t.Cell(1, 5)
$old = $cell.Range.Text
if ($old.Substring(0, $old.Length - 2) -ne '56÷3=18, 2') {
    throw "Cell(1,5) expected '56÷3=18, 2' but found '$old'"
}
$cell.Range.Text = '53÷9=5, 8'

$cell = $t.Cell(5, 1)
$old = $cell.Range.Text
if ($old.Substring(0, $old.Length - 2) -ne '70÷4=17, 2') {
    throw "Cell(5,1) expected '70÷4=17, 2' but found '$old'"
}
$cell.Range.Text = '36÷4=9, 0'

$cell = $t.Cell(5, 2)
$old = $cell.Range.Text
if ($old.Substring(0, $old.Length - 2) -ne '48÷6=8, 0') {
    throw "Cell(5,2) expected '48÷6=8, 0' but found '$old'"
}
$cell.Range.Text = '25÷6=4, 1'

$cell = $t.Cell(5, 3)
$old = $cell.Range.Text
if ($old.Substring(0, $old.Length - 2) -ne '93÷5=18, 3') {
    throw "Cell(5,3) expected '93÷5=18, 3' but found '$old'"
}
$cell.Range.Text = '76÷3=25, 1'

$cell = $t.Cell(5, 4)
$old = $cell.Range.Text
if ($old.Substring(0, $old.Length - 2) -ne '62÷9=6, 8') {
    throw "Cell(5,4) expected '62÷9=6, 8' but found '$old'"
}
$cell.Range.Text = '22÷8=2, 6'

$cell = $t.Cell(5, 5)
$old = $cell.Range.Text
if ($old.Substring(0, $old.Length - 2) -ne '14÷6=2, 2') {
    throw "Cell(5,5) expected '14÷6=2, 2' but found '$old'"
}
$cell.Range.Text = '10÷9=1, 1'

$cell = $t.Cell(9, 1)
$old = $cell.Range.Text
if ($old.Substring(0, $old.Length - 2) -ne '79÷8=9, 7') {
    throw "Cell(9,1) expected '79÷8=9, 7' but found '$old'"
}
$cell.Range.Text = '51÷9=5, 6'

$cell = $t.Cell(9, 2)
$old = $cell.Range.Text
if ($old.Substring(0, $old.Length - 2) -ne '27÷8=3, 3') {
    throw "Cell(9,2) expected '27÷8=3, 3' but found '$old'"
}
$cell.Range.Text = '13÷2=6, 1'

$cell = $t.Cell(9, 3)
$old = $cell.Range.Text
if ($old.Substring(0, $old.Length - 2) -ne '26÷2=13, 0') {
    throw "Cell(9,3) expected '26÷2=13, 0' but found '$old'"
}
$cell.Range.Text = '37÷2=18, 1'

$cell = $t.Cell(9, 4)
$old = $cell.Range.Text
if ($old.Substring(0, $old.Length - 2) -ne '36÷3=12, 0') {
    throw "Cell(9,4) expected '36÷3=12, 0' but found '$old'"
}
$cell.Range.Text = '96÷9=10, 6'

$cell = $t.Cell(9, 5)
$old = $cell.Range.Text
if ($old.Substring(0, $old.Length - 2) -ne '47÷3=15, 2') {
    throw "Cell(9,5) expected '47÷3=15, 2' but found '$old'"
}
$cell.Range.Text = '76÷2=38, 0'

$cell = $t.Cell(13, 1)
$old = $cell.Range.Text
if ($old.Substring(0, $old.Length - 2) -ne '94÷5=18, 4') {
    throw "Cell(13,1) expected '94÷5=18, 4' but found '$old'"
}
$cell.Range.Text = '99÷3=33, 0'

$cell = $t.Cell(13, 2)
$old = $cell.Range.Text
if ($old.Substring(0, $old.Length - 2) -ne '20÷9=2, 2') {
    throw "Cell(13,2) expected '20÷9=2, 2' but found '$old'"
}
$cell.Range.Text = '25÷9=2, 7'

$cell = $t.Cell(13, 3)
$old = $cell.Range.Text
if ($old.Substring(0, $old.Length - 2) -ne '91÷4=22, 3') {
    throw "Cell(13,3) expected '91÷4=22, 3' but found '$old'"
}
$cell.Range.Text = '29÷6=4, 5'

$cell = $t.Cell(13, 4)
$old = $cell.Range.Text
if ($old.Substring(0, $old.Length - 2) -ne '40÷9=4, 4') {
    throw "Cell(13,4) expected '40÷9=4, 4' but found '$old'"
}
$cell.Range.Text = '96÷4=24, 0'

$cell = $t.Cell(13, 5)
$old = $cell.Range.Text
if ($old.Substring(0, $old.Length - 2) -ne '96÷9=10, 6') {
    throw "Cell(13,5) expected '96÷9=10, 6' but found '$old'"
}
$cell.Range.Text = '90÷6=15, 0'

$cell = $t.Cell(17, 1)
$old = $cell.Range.Text
if ($old.Substring(0, $old.Length - 2) -ne '25÷5=5, 0') {
    throw "Cell(17,1) expected '25÷5=5, 0' but found '$old'"
}
$cell.Range.Text = '87÷2=43, 1'

$cell = $t.Cell(17, 2)
$old = $cell.Range.Text
if ($old.Substring(0, $old.Length - 2) -ne '35÷3=11, 2') {
    throw "Cell(17,2) expected '35÷3=11, 2' but found '$old'"
}
$cell.Range.Text = '57÷9=6, 3'

$cell = $t.Cell(17, 3)
$old = $cell.Range.Text
if ($old.Substring(0, $old.Length - 2) -ne '61÷9=6, 7') {
    throw "Cell(17,3) expected '61÷9=6, 7' but found '$old'"
}
$cell.Range.Text = '42÷4=10, 2'

$cell = $t.Cell(17, 4)
$old = $cell.Range.Text
if ($old.Substring(0, $old.Length - 2) -ne '48÷5=9, 3') {
    throw "Cell(17,4) expected '48÷5=9, 3' but found '$old'"
}
$cell.Range.Text = '57÷7=8, 1'

$cell = $t.Cell(17, 5)
$old = $cell.Range.Text
if ($old.Substring(0, $old.Length - 2) -ne '88÷3=29, 1') {
    throw "Cell(17,5) expected '88÷3=29, 1' but found '$old'"
}
$cell.Range.Text = '49÷3=16, 1'
